# Rename the inline logo pictures that live in the document's headers and
# footers. The Pearson logo (alt text/description points at the
# "PearsonLogo.png" asset) goes from "image2.png" -> "image1.png", and the
# BTec logo (alt text/description "BTec_Logo-Orange") goes from
# "image1.jpg" -> "image2.jpg". Both logos appear twice (the "first page"
# header/footer and the "default" header/footer), so walk every
# header/footer of every section and rename by matching on the picture's
# AlternativeText/description, which uniquely identifies each logo.

$d = $word.ActiveDocument

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($h = 1; $h -le 3; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            $r = $hdr.Range
            for ($i = 1; $i -le $r.InlineShapes.Count; $i++) {
                $shp = $r.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                } elseif ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }

        $ftr = $sec.Footers.Item($h)
        if ($ftr.Exists) {
            $r2 = $ftr.Range
            for ($i = 1; $i -le $r2.InlineShapes.Count; $i++) {
                $shp2 = $r2.InlineShapes.Item($i)
                if ($shp2.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp2.Name = "image2.jpg"
                } elseif ($shp2.AlternativeText -like "*PearsonLogo.png") {
                    $shp2.Name = "image1.png"
                }
            }
        }
    }
}
